$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2.565309470876655
$ws.Range("D3").Value = 3.003081339182493
$ws.Range("D4").Value = 4.199891207049152
$ws.Range("D5").Value = 3.274189803507521
$ws.Range("D6").Value = 8.725130363174229
$ws.Range("D7").Value = 4.720781948531094
$ws.Range("D8").Value = 2.535127651661789
$ws.Range("D9").Value = 4.809824706622224
$ws.Range("D10").Value = 6.13938237555457
$ws.Range("D11").Value = 2.86069802887924
$ws.Range("D12").Value = 2.30361538158524
$ws.Range("D13").Value = 6.547923009847037
$ws.Range("D14").Value = 4.928269433743423
$ws.Range("D15").Value = 9.531842630873502
$ws.Range("D16").Value = 4.323533340851341
$ws.Range("D17").Value = 12.51502503806533
$ws.Range("D18").Value = 3.525476847061302
$ws.Range("D19").Value = 25.37179967584758
$ws.Range("D20").Value = 10.36424620600918
$ws.Range("D21").Value = 9.435240225682612
$ws.Range("D22").Value = 16.09327796865148
$ws.Range("D23").Value = 18.30366463356607
$ws.Range("D24").Value = 6.111661780494481
$ws.Range("D25").Value = 16.64950296625551
$ws.Range("D26").Value = 24.8828281198701
$ws.Range("D27").Value = 16.1454417455465
